$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "gender" column (I) - the gender-mapping feature added for issue #48.
# Written in "gender, m, male, female" first-use order so the shared-string
# table is built up the same way it was in the authored workbook.
$ws.Range("I1").Value = "gender"
$ws.Range("I3").Value = "m"
$ws.Range("I5").Value = "male"
$ws.Range("I2").Value = "female"
$ws.Range("I4").Value = "m"
$ws.Range("I6").Value = "m"
$ws.Range("I7").Value = "m"
$ws.Range("I8").Value = "male"

$ws.Range("I1:I8").HorizontalAlignment = -4108

# Matches the saved file's active selection after the edit.
$ws.Range("I9").Select()
